# Apply math-expression updates to the "within100" practice-sheet table.
# Each cell in the 20x5 table holds a single arithmetic expression such as
# "63+3=". This updates them in place, cell by cell, preserving the
# surrounding run/paragraph formatting (font, size, etc.) already present
# in each cell.

function Set-CellText {
    param($table, $row, $col, $expectedOld, $newText)

    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range

    # Cell.Range.Text includes the trailing cell-mark / paragraph-mark
    # characters (chr 13 + chr 7); strip them before comparing so we can
    # sanity-check we are editing the expected cell.
    $actual = $cellRange.Text
    $actualTrimmed = $actual.TrimEnd([char]7, [char]13)
    if ($actualTrimmed -ne $expectedOld) {
        Write-Host "WARNING: cell ($row,$col) expected '$expectedOld' but found '$actualTrimmed'"
    }

    $cellRange.Text = $newText
}

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

Set-CellText $tbl 1 1 "63+3=" "0+14="
Set-CellText $tbl 1 2 "77-15=" "72+7="
Set-CellText $tbl 1 3 "43-15=" "31+34="
Set-CellText $tbl 1 4 "88-29=" "19+56="
Set-CellText $tbl 1 5 "53+10=" "98-88="
Set-CellText $tbl 2 1 "22+0=" "19+12="
Set-CellText $tbl 2 2 "12+67=" "29-21="
Set-CellText $tbl 2 3 "49+33=" "17+40="
Set-CellText $tbl 2 4 "12+70=" "30+53="
Set-CellText $tbl 2 5 "55-43=" "76+12="
Set-CellText $tbl 3 1 "59-25=" "38-17="
Set-CellText $tbl 3 2 "7+26=" "59+21="
Set-CellText $tbl 3 3 "21+50=" "81-72="
Set-CellText $tbl 3 4 "24-18=" "6+70="
Set-CellText $tbl 3 5 "99-21=" "91-89="
Set-CellText $tbl 4 1 "83-43=" "14+20="
Set-CellText $tbl 4 2 "74-19=" "82-24="
Set-CellText $tbl 4 3 "69-9=" "90-41="
Set-CellText $tbl 4 4 "66+29=" "39+48="
Set-CellText $tbl 4 5 "49-45=" "77-6="
Set-CellText $tbl 5 1 "54-47=" "44-4="
Set-CellText $tbl 5 2 "23+52=" "14+67="
Set-CellText $tbl 5 3 "62-28=" "98-51="
Set-CellText $tbl 5 4 "27-10=" "57+25="
Set-CellText $tbl 5 5 "29-5=" "24+72="
Set-CellText $tbl 6 1 "76-20=" "92-61="
Set-CellText $tbl 6 2 "74+24=" "59+40="
Set-CellText $tbl 6 3 "81-7=" "2+35="
Set-CellText $tbl 6 4 "49+34=" "58-30="
Set-CellText $tbl 6 5 "61-27=" "31+39="
Set-CellText $tbl 7 1 "94-32=" "51-30="
Set-CellText $tbl 7 2 "21+16=" "57+7="
Set-CellText $tbl 7 3 "53-13=" "4+80="
Set-CellText $tbl 7 4 "8+5=" "41-2="
Set-CellText $tbl 7 5 "7-5=" "64+1="
Set-CellText $tbl 8 1 "90-49=" "40+34="
Set-CellText $tbl 8 2 "26-17=" "58-57="
Set-CellText $tbl 8 3 "81+9=" "84-61="
Set-CellText $tbl 8 4 "45+32=" "9+29="
Set-CellText $tbl 8 5 "3+19=" "67-67="
Set-CellText $tbl 9 1 "24+66=" "29-13="
Set-CellText $tbl 9 2 "55-43=" "82-52="
Set-CellText $tbl 9 3 "84-6=" "99-71="
Set-CellText $tbl 9 4 "58+8=" "29+46="
Set-CellText $tbl 9 5 "26-12=" "74+9="
Set-CellText $tbl 10 1 "70-66=" "94-59="
Set-CellText $tbl 10 2 "23+40=" "82-16="
Set-CellText $tbl 10 3 "25+1=" "69-33="
Set-CellText $tbl 10 4 "84-5=" "72+7="
Set-CellText $tbl 10 5 "73+21=" "94+3="
Set-CellText $tbl 11 1 "74-27=" "27+1="
Set-CellText $tbl 11 2 "19+64=" "29+52="
Set-CellText $tbl 11 3 "90-28=" "47+36="
Set-CellText $tbl 11 4 "80-64=" "43+14="
Set-CellText $tbl 11 5 "93-29=" "75+2="
Set-CellText $tbl 12 1 "6+86=" "36+15="
Set-CellText $tbl 12 2 "24-15=" "1+97="
Set-CellText $tbl 12 3 "17+18=" "5+66="
Set-CellText $tbl 12 4 "42+6=" "3+38="
Set-CellText $tbl 12 5 "25-22=" "95-51="
Set-CellText $tbl 13 1 "39+0=" "29-26="
Set-CellText $tbl 13 2 "40-13=" "49-42="
Set-CellText $tbl 13 3 "77+12=" "26+66="
Set-CellText $tbl 13 4 "7+88=" "43+11="
Set-CellText $tbl 13 5 "82-39=" "75-20="
Set-CellText $tbl 14 1 "28+9=" "73-62="
Set-CellText $tbl 14 2 "35+57=" "27+54="
Set-CellText $tbl 14 3 "47-17=" "74-6="
Set-CellText $tbl 14 4 "12+77=" "14+40="
Set-CellText $tbl 14 5 "83-65=" "62-5="
Set-CellText $tbl 15 1 "49-5=" "33+56="
Set-CellText $tbl 15 2 "89+0=" "54-51="
Set-CellText $tbl 15 3 "20+69=" "60-34="
Set-CellText $tbl 15 4 "76-73=" "62-41="
Set-CellText $tbl 15 5 "66-9=" "61-15="
Set-CellText $tbl 16 1 "20+17=" "33+14="
Set-CellText $tbl 16 2 "33+21=" "97-55="
Set-CellText $tbl 16 3 "4+10=" "3+65="
Set-CellText $tbl 16 4 "64+5=" "48-7="
Set-CellText $tbl 16 5 "50+31=" "38+3="
Set-CellText $tbl 17 1 "45+44=" "47+21="
Set-CellText $tbl 17 2 "1+56=" "44+16="
Set-CellText $tbl 17 3 "4+31=" "30-3="
Set-CellText $tbl 17 4 "80-46=" "26+56="
Set-CellText $tbl 17 5 "17-16=" "91-78="
Set-CellText $tbl 18 1 "8+77=" "19+10="
Set-CellText $tbl 18 2 "47+45=" "37+51="
Set-CellText $tbl 18 3 "21-11=" "73-3="
Set-CellText $tbl 18 4 "70-66=" "51+21="
Set-CellText $tbl 18 5 "88-73=" "8+71="
Set-CellText $tbl 19 1 "51+17=" "41+16="
Set-CellText $tbl 19 2 "71-13=" "46-34="
Set-CellText $tbl 19 3 "36-18=" "1+46="
Set-CellText $tbl 19 4 "35+57=" "26+6="
Set-CellText $tbl 19 5 "9+7=" "65-19="
Set-CellText $tbl 20 1 "51+22=" "17+3="
Set-CellText $tbl 20 2 "84-63=" "25-3="
Set-CellText $tbl 20 3 "79-64=" "43+3="
Set-CellText $tbl 20 4 "49+22=" "85-29="
Set-CellText $tbl 20 5 "61+2=" "92-59="
